$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert three new rows above row 102 (old row 102 "If difficult..." shifts to row 105) ---
$ws.Rows.Item(102).Insert()
$ws.Rows.Item(102).Insert()
$ws.Rows.Item(102).Insert()

# --- New row 102: "Cable harnesses:" header + first AliExpress link (plain text) ---
$ws.Cells.Item(102, 1).HorizontalAlignment = -4108
$ws.Cells.Item(102, 1).VerticalAlignment = -4108

$ws.Cells.Item(102, 2).Value = "Cable harnesses:"
$ws.Cells.Item(102, 2).NumberFormat = "@"
$ws.Cells.Item(102, 2).HorizontalAlignment = -4131
$ws.Cells.Item(102, 2).VerticalAlignment = -4108
$ws.Cells.Item(102, 2).IndentLevel = 1

$ws.Cells.Item(102, 3).Value = "https://nl.aliexpress.com/item/1005002942389730.html?spm=a2g0o.productlist.main.9.8badZ4V9Z4V92a&algo_pvid=94a6380a-288b-4c3a-acc3-61c2ba7c95ce&algo_exp_id=94a6380a-288b-4c3a-acc3-61c2ba7c95ce-8&pdp_ext_f=%7B%22order%22%3A%2241%22%2C%22eval%22%3A%221%22%7D&pdp_npi=6%40dis%21EUR%210.31%210.27%21%21%210.35%210.30%21%40211b80d117550731005353576e8de1%2112000022896959057%21sea%21NL%210%21ABX%211%210%21n_tag%3A-29910%3Bm03_new_user%3A-29895&curPageLogUid=Souz7hgWNLUp&utparam-url=scene%3Asearch%7Cquery_from%3A%7Cx_object_id%3A1005002942389730%7C_p_origin_prod%3A"
$ws.Cells.Item(102, 3).HorizontalAlignment = -4131
$ws.Cells.Item(102, 3).VerticalAlignment = -4108
$ws.Cells.Item(102, 3).IndentLevel = 1

for ($col = 4; $col -le 11; $col++) {
  $ws.Cells.Item(102, $col).HorizontalAlignment = -4131
  $ws.Cells.Item(102, $col).VerticalAlignment = -4108
  $ws.Cells.Item(102, $col).IndentLevel = 1
}

# --- New row 103: second AliExpress link, turned into a real hyperlink ---
# (Hyperlinks.Add first, while C103 is still unformatted, so Excel only swaps
# in the Hyperlink named style without carrying over any other attributes.)
$ws.Hyperlinks.Add($ws.Cells.Item(103, 3), "https://nl.aliexpress.com/w/wholesale-pre-crimped-cable-6-colors.html?spm=a2g0o.home.auto_suggest.4.60ca306bMkCXSW", "", "", " https://nl.aliexpress.com/w/wholesale-pre-crimped-cable-6-colors.html?spm=a2g0o.home.auto_suggest.4.60ca306bMkCXSW")

$ws.Cells.Item(103, 1).HorizontalAlignment = -4108
$ws.Cells.Item(103, 1).VerticalAlignment = -4108

$ws.Cells.Item(103, 2).NumberFormat = "@"
$ws.Cells.Item(103, 2).HorizontalAlignment = -4131
$ws.Cells.Item(103, 2).VerticalAlignment = -4108
$ws.Cells.Item(103, 2).IndentLevel = 1

for ($col = 4; $col -le 11; $col++) {
  $ws.Cells.Item(103, $col).HorizontalAlignment = -4131
  $ws.Cells.Item(103, $col).VerticalAlignment = -4108
  $ws.Cells.Item(103, $col).IndentLevel = 1
}

# --- Update selection / scroll position ---
$ws.Range("F112").Select()
